$wb = $excel.ActiveWorkbook

# --- Netherlands (copy of Greece) ---
$greece = $wb.Worksheets.Item("Greece")
$greece.Copy($null, $greece)
$nl = $wb.Worksheets.Item("Greece (2)")
$nl.Name = "Netherlands"
$nl.Range("B4").Value = "NGC-3144/T2199"
$nl.Range("B2").Value = "Netherlands Market"
$a9 = $nl.Range("A9").Value2
$a10 = $nl.Range("A10").Value2
$nl.Range("A9").Value = $a10
$nl.Range("A10").Value = $a9
$nl.Range("B1:D12").EntireColumn.AutoFit() | Out-Null

# --- Austria (copy of Netherlands, PR1D2-Unmonitored row removed) ---
$nl.Copy($null, $nl)
$at = $wb.Worksheets.Item("Netherlands (2)")
$at.Name = "Austria"
$at.Range("B4").Value = "NGC-3817/T2306"
$at.Range("B2").Value = "Austria Market"
$at.Range("A10").EntireRow.Delete() | Out-Null
$at.Range("B1:D11").EntireColumn.AutoFit() | Out-Null

# --- Denmark (copy of Austria) ---
$at.Copy($null, $at)
$dk = $wb.Worksheets.Item("Austria (2)")
$dk.Name = "Denmark"
$dk.Range("B4").Value = "NGC-2913/T2306"
$dk.Range("B2").Value = "Denmark Market"
$dk.Range("B1:D11").EntireColumn.AutoFit() | Out-Null

# Netherlands ends up the active/selected sheet
$nl.Activate()
